# Updates cryptos list values (Price and Volume(1h) columns)
# Leading apostrophe forces text interpretation (matching original inlineStr cell type)
# without leaving stray numeric values or altering number formats.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.890.29"
$ws.Range("E2").Value = "'  -0.24%  "

$ws.Range("D3").Value = "'1.584.16"
$ws.Range("E3").Value = "'  -2.02%  "

$ws.Range("E4").Value = "'  -0.27%  "

$ws.Range("D5").Value = "'210.24"
$ws.Range("E5").Value = "'  -0.62%  "

$ws.Range("E7").Value = "'  -2.13%  "

$ws.Range("E8").Value = "'  +0.16%  "

$ws.Range("D9").Value = "'0.0613"
$ws.Range("E9").Value = "'  -1.37%  "

$ws.Range("D10").Value = "'18.08"
$ws.Range("E10").Value = "'  -0.15%  "

$ws.Range("E11").Value = "'  -0.11%  "

$ws.Range("D12").Value = "'1.804.98"
$ws.Range("E12").Value = "'  -1.93%  "

$ws.Range("D13").Value = "'1.581.56"
$ws.Range("E13").Value = "'  -2.22%  "

$ws.Range("E14").Value = "'  -2.23%  "

$ws.Range("E15").Value = "'  -2.07%  "

$ws.Range("D16").Value = "'25.877.44"
$ws.Range("E16").Value = "'  -0.28%  "

$ws.Range("E17").Value = "'  -0.70%  "

$ws.Range("D18").Value = "'59.97"
$ws.Range("E18").Value = "'  -2.44%  "

$ws.Range("E19").Value = "'  -0.25%  "

$ws.Range("D20").Value = "'192.96"
$ws.Range("E20").Value = "'  +1.07%  "

$ws.Range("E21").Value = "'  -0.54%  "

$ws.Range("D22").Value = "'9.35"
$ws.Range("E22").Value = "'  -0.51%  "

$ws.Range("D23").Value = "'5.93"
$ws.Range("E23").Value = "'  -0.94%  "

$ws.Range("E24").Value = "'  -0.20%  "

$ws.Range("D25").Value = "'141.65"
$ws.Range("E25").Value = "'  -0.98%  "

$ws.Range("E26").Value = "'  -0.34%  "

$ws.Range("D27").Value = "'1.70"
$ws.Range("E27").Value = "'  -0.40%  "

$ws.Range("E28").Value = "'  -0.28%  "

$ws.Range("E29").Value = "'  -2.33%  "

$ws.Range("E30").Value = "'  -4.77%  "

$ws.Range("E31").Value = "'  -0.30%  "

$ws.Range("E32").Value = "'  +0.31%  "

$ws.Range("E33").Value = "'  -1.44%  "

$ws.Range("E34").Value = "'  +0.88%  "

$ws.Range("E35").Value = "'  -2.19%  "

$ws.Range("D36").Value = "'1.096.57"
$ws.Range("E36").Value = "'  -2.24%  "

$ws.Range("E37").Value = "'  -0.38%  "

$ws.Range("E38").Value = "'  -1.87%  "

$ws.Range("D39").Value = "'0.0151"
$ws.Range("E39").Value = "'  -0.62%  "

$ws.Range("E40").Value = "'  -2.39%  "

$ws.Range("E41").Value = "'  -4.20%  "

$ws.Range("D42").Value = "'0.799"
$ws.Range("E42").Value = "'  +5.85%  "

$ws.Range("D43").Value = "'93.43"

$ws.Range("D44").Value = "'5.12"
$ws.Range("E44").Value = "'  +1.32%  "

$ws.Range("D45").Value = "'1.718.72"
$ws.Range("E45").Value = "'  -1.87%  "

$ws.Range("E46").Value = "'  -1.65%  "

$ws.Range("E47").Value = "'  +2.15%  "

$ws.Range("D48").Value = "'53.12"
$ws.Range("E48").Value = "'  -1.01%  "

$ws.Range("E49").Value = "'  -1.26%  "

$ws.Range("E50").Value = "'  -1.05%  "

$ws.Range("E51").Value = "'  -0.16%  "
